$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2579443333333333
$ws.Range("H2").Value = 0.773833
$ws.Range("I2").Value = 0.05467096027587241
$ws.Range("J2").Value = 0.0546709602758724
$ws.Range("M2").Value = 8.415202000000001
$ws.Range("N2").Value = 25.245606
$ws.Range("O2").Value = 0.1569653516800918
$ws.Range("P2").Value = 0.1569653516800918
$ws.Range("Q2").Value = 2.170653669755334
$ws.Range("R2").Value = 19.535883027798
$ws.Range("S2").Value = 0.008581446506390639
$ws.Range("T2").Value = 0.008581446506390637

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2579443333333333
$ws.Range("H3").Value = 0.773833
$ws.Range("I3").Value = 0.05467096027587241
$ws.Range("J3").Value = 0.0546709602758724
$ws.Range("O3").Value = 0.5328513631375226
$ws.Range("P3").Value = 0.5328513631375226
$ws.Range("Q3").Value = 7.368733000298778
$ws.Range("R3").Value = 66.31859700268899
$ws.Range("S3").Value = 0.02913149570703596
$ws.Range("T3").Value = 0.02913149570703595

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2579443333333333
$ws.Range("H4").Value = 0.773833
$ws.Range("I4").Value = 0.05467096027587241
$ws.Range("J4").Value = 0.0546709602758724
$ws.Range("M4").Value = 2.790681000000001
$ws.Range("N4").Value = 8.372043000000001
$ws.Range("O4").Value = 0.05205344144940909
$ws.Range("P4").Value = 0.05205344144940908
$ws.Range("Q4").Value = 0.7198403500910001
$ws.Range("R4").Value = 6.478563150819001
$ws.Range("S4").Value = 0.002845811629703094
$ws.Range("T4").Value = 0.002845811629703094

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2579443333333333
$ws.Range("H5").Value = 0.773833
$ws.Range("I5").Value = 0.05467096027587241
$ws.Range("J5").Value = 0.0546709602758724
$ws.Range("M5").Value = 13.838817
$ws.Range("N5").Value = 41.516451
$ws.Range("O5").Value = 0.2581298437329766
$ws.Range("P5").Value = 0.2581298437329766
$ws.Range("Q5").Value = 3.569644425187
$ws.Range("R5").Value = 32.126799826683
$ws.Range("S5").Value = 0.01411220643274272
$ws.Range("T5").Value = 0.01411220643274271

$ws.Range("I6").Value = 0.4403936734732808
$ws.Range("J6").Value = 0.4403936734732807
$ws.Range("M6").Value = 8.415202000000001
$ws.Range("N6").Value = 25.245606
$ws.Range("O6").Value = 0.1569653516800918
$ws.Range("P6").Value = 0.1569653516800918
$ws.Range("Q6").Value = 17.48537319699667
$ws.Range("R6").Value = 157.36835877297
$ws.Range("S6").Value = 0.06912654783442101
$ws.Range("T6").Value = 0.06912654783442101

$ws.Range("I7").Value = 0.4403936734732808
$ws.Range("J7").Value = 0.4403936734732807
$ws.Range("O7").Value = 0.5328513631375226
$ws.Range("P7").Value = 0.5328513631375226
$ws.Range("S7").Value = 0.2346643692273787
$ws.Range("T7").Value = 0.2346643692273787

$ws.Range("I8").Value = 0.4403936734732808
$ws.Range("J8").Value = 0.4403936734732807
$ws.Range("M8").Value = 2.790681000000001
$ws.Range("N8").Value = 8.372043000000001
$ws.Range("O8").Value = 0.05205344144940909
$ws.Range("P8").Value = 0.05205344144940908
$ws.Range("Q8").Value = 5.798565353365002
$ws.Range("R8").Value = 52.18708818028502
$ws.Range("S8").Value = 0.0229240062968316
$ws.Range("T8").Value = 0.0229240062968316

$ws.Range("I9").Value = 0.4403936734732808
$ws.Range("J9").Value = 0.4403936734732807
$ws.Range("M9").Value = 13.838817
$ws.Range("N9").Value = 41.516451
$ws.Range("O9").Value = 0.2581298437329766
$ws.Range("P9").Value = 0.2581298437329766
$ws.Range("Q9").Value = 28.754732191805
$ws.Range("R9").Value = 258.7925897262451
$ws.Range("S9").Value = 0.1136787501146495
$ws.Range("T9").Value = 0.1136787501146495

$ws.Range("G10").Value = 2.217259
$ws.Range("H10").Value = 6.651777
$ws.Range("I10").Value = 0.4699451123575263
$ws.Range("J10").Value = 0.4699451123575263
$ws.Range("M10").Value = 8.415202000000001
$ws.Range("N10").Value = 25.245606
$ws.Range("O10").Value = 0.1569653516800918
$ws.Range("P10").Value = 0.1569653516800918
$ws.Range("Q10").Value = 18.658682371318
$ws.Range("R10").Value = 167.928141341862
$ws.Range("S10").Value = 0.07376509983153937
$ws.Range("T10").Value = 0.07376509983153937

$ws.Range("G11").Value = 2.217259
$ws.Range("H11").Value = 6.651777
$ws.Range("I11").Value = 0.4699451123575263
$ws.Range("J11").Value = 0.4699451123575263
$ws.Range("O11").Value = 0.5328513631375226
$ws.Range("P11").Value = 0.5328513631375226
$ws.Range("Q11").Value = 63.34075787738232
$ws.Range("R11").Value = 570.066820896441
$ws.Range("S11").Value = 0.2504108937195241
$ws.Range("T11").Value = 0.2504108937195241

$ws.Range("G12").Value = 2.217259
$ws.Range("H12").Value = 6.651777
$ws.Range("I12").Value = 0.4699451123575263
$ws.Range("J12").Value = 0.4699451123575263
$ws.Range("M12").Value = 2.790681000000001
$ws.Range("N12").Value = 8.372043000000001
$ws.Range("O12").Value = 0.05205344144940909
$ws.Range("P12").Value = 0.05205344144940908
$ws.Range("Q12").Value = 6.187662563379001
$ws.Range("R12").Value = 55.68896307041101
$ws.Range("S12").Value = 0.02446226039053847
$ws.Range("T12").Value = 0.02446226039053847

$ws.Range("G13").Value = 2.217259
$ws.Range("H13").Value = 6.651777
$ws.Range("I13").Value = 0.4699451123575263
$ws.Range("J13").Value = 0.4699451123575263
$ws.Range("M13").Value = 13.838817
$ws.Range("N13").Value = 41.516451
$ws.Range("O13").Value = 0.2581298437329766
$ws.Range("P13").Value = 0.2581298437329766
$ws.Range("Q13").Value = 30.684241542603
$ws.Range("R13").Value = 276.158173883427
$ws.Range("S13").Value = 0.1213068584159244
$ws.Range("T13").Value = 0.1213068584159244

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1650883333333333
$ws.Range("H14").Value = 0.495265
$ws.Range("I14").Value = 0.03499025389332058
$ws.Range("J14").Value = 0.03499025389332058
$ws.Range("M14").Value = 8.415202000000001
$ws.Range("N14").Value = 25.245606
$ws.Range("O14").Value = 0.1569653516800918
$ws.Range("P14").Value = 0.1569653516800918
$ws.Range("Q14").Value = 1.389251672843334
$ws.Range("R14").Value = 12.50326505559
$ws.Range("S14").Value = 0.005492257507740766
$ws.Range("T14").Value = 0.005492257507740764

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1650883333333333
$ws.Range("H15").Value = 0.495265
$ws.Range("I15").Value = 0.03499025389332058
$ws.Range("J15").Value = 0.03499025389332058
$ws.Range("O15").Value = 0.5328513631375226
$ws.Range("P15").Value = 0.5328513631375226
$ws.Range("Q15").Value = 4.716102246082778
$ws.Range("R15").Value = 42.444920214745
$ws.Range("S15").Value = 0.01864460448358388
$ws.Range("T15").Value = 0.01864460448358388

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1650883333333333
$ws.Range("H16").Value = 0.495265
$ws.Range("I16").Value = 0.03499025389332058
$ws.Range("J16").Value = 0.03499025389332058
$ws.Range("M16").Value = 2.790681000000001
$ws.Range("N16").Value = 8.372043000000001
$ws.Range("O16").Value = 0.05205344144940909
$ws.Range("P16").Value = 0.05205344144940908
$ws.Range("Q16").Value = 0.4607088751550001
$ws.Range("R16").Value = 4.146379876395001
$ws.Range("S16").Value = 0.001821363132335921
$ws.Range("T16").Value = 0.001821363132335921

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1650883333333333
$ws.Range("H17").Value = 0.495265
$ws.Range("I17").Value = 0.03499025389332058
$ws.Range("J17").Value = 0.03499025389332058
$ws.Range("M17").Value = 13.838817
$ws.Range("N17").Value = 41.516451
$ws.Range("O17").Value = 0.2581298437329766
$ws.Range("P17").Value = 0.2581298437329766
$ws.Range("Q17").Value = 2.284627233835
$ws.Range("R17").Value = 20.561645104515
$ws.Range("S17").Value = 0.009032028769660017
$ws.Range("T17").Value = 0.009032028769660016
